$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 6327
$ws1.Range("F4").Value = 178
$ws1.Range("F7").Value = 1903
$ws1.Range("F8").Value = 1439
$ws1.Range("F9").Value = 297
$ws1.Range("F11").Value = 266
$ws1.Range("F12").Value = 5588

# Sheet "全部类型" (fourth sheet) - same rows, same metric, slightly different F3 value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 6328
$ws4.Range("F4").Value = 178
$ws4.Range("F7").Value = 1903
$ws4.Range("F8").Value = 1439
$ws4.Range("F9").Value = 297
$ws4.Range("F11").Value = 266
$ws4.Range("F12").Value = 5588
